$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Year of Treatment" column (column B). Everything to the
# right (Alone, With the family of origin..., ..., Total) shifts left
# by one column.
$ws.Range("B1").EntireColumn.Delete()
